# Insert two new price-report rows (Lane Late) right before the existing
# row 470, shifting all subsequent rows down by two (matches the dimension
# growing from A1:T568 to A1:T570).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("470:471").Insert()

# New row 470
$ws.Range("A470").Value = 3
$ws.Range("B470").Value = "Femacal de La Calera"
$ws.Range("C470").Value = "Coquimbo"
$ws.Range("D470").Value = 44543
$ws.Range("E470").Value = 5
$ws.Range("F470").Value = "Fruta"
$ws.Range("G470").Value = 100102
$ws.Range("H470").Value = "Cítricos"
$ws.Range("I470").Value = 100102005
$ws.Range("J470").Value = "Naranja"
$ws.Range("K470").Value = "Lane Late"
$ws.Range("L470").Value = "Primera"
$ws.Range("M470").Value = 56
$ws.Range("N470").Value = 8000
$ws.Range("O470").Value = 8000
$ws.Range("P470").Value = 8000
$ws.Range("Q470").Value = '$/malla 13 kilos'
$ws.Range("R470").Value = "Provincia de Quillota"
$ws.Range("S470").Value = 615
$ws.Range("T470").Value = 13

# New row 471
$ws.Range("A471").Value = 3
$ws.Range("B471").Value = "Femacal de La Calera"
$ws.Range("C471").Value = "Coquimbo"
$ws.Range("D471").Value = 44543
$ws.Range("E471").Value = 5
$ws.Range("F471").Value = "Fruta"
$ws.Range("G471").Value = 100102
$ws.Range("H471").Value = "Cítricos"
$ws.Range("I471").Value = 100102005
$ws.Range("J471").Value = "Naranja"
$ws.Range("K471").Value = "Lane Late"
$ws.Range("L471").Value = "Segunda"
$ws.Range("M471").Value = 60
$ws.Range("N471").Value = 6000
$ws.Range("O471").Value = 6000
$ws.Range("P471").Value = 6000
$ws.Range("Q471").Value = '$/malla 13 kilos'
$ws.Range("R471").Value = "Provincia de Quillota"
$ws.Range("S471").Value = 462
$ws.Range("T471").Value = 13
